# "C - Decision control started"
# Introduce an `if(a[j]<a[i]) { ... }` decision block inside the nested
# for(j...) loop of the selection-sort routine.

$d = $word.ActiveDocument

# --- Step 1: the old "TAB TAB for(j= i + 1 ; j < n ; j++)" paragraph loses
#     its two leading tabs and instead gets hanging-indent formatting.
$pFor = $d.Paragraphs.Item(35)
$rFor = $d.Range($pFor.Range.Start, $pFor.Range.End - 1)
$rFor.Text = "for(j= i + 1 ; j < n ; j++)"
$pFor.LeftIndent = 72
$pFor.FirstLineIndent = 36

# --- Step 2: new paragraph carrying the opening brace of the for(j) body.
$pFor.Range.InsertParagraphAfter()
$pBrace1 = $d.Paragraphs.Item(36)
$pBrace1.LeftIndent = 36
$pBrace1.FirstLineIndent = 0
$pBrace1.Range.InsertAfter("`t`t{")

# --- Step 3: new paragraph for the "if(a[j]<a[i])" condition.
$pBrace1.Range.InsertParagraphAfter()
$pIf = $d.Paragraphs.Item(37)
$pIf.LeftIndent = 108
$pIf.FirstLineIndent = 36
$pIf.Range.InsertAfter("if(a[j]<a[i])")

# --- Step 4: new paragraph for the opening brace of the if body (this is
#     where Word's last-edit "_GoBack" marker will now sit).
$pIf.Range.InsertParagraphAfter()
$pBrace2 = $d.Paragraphs.Item(38)
$pBrace2.LeftIndent = 36
$pBrace2.FirstLineIndent = 0
$pBrace2.Range.InsertAfter("`t`t`t{")

Write-Output "done step 1-4"
